# Proper handling of STAGE with data files (ignore if datarecord stage <> stage from globals)
# Adds duplicate "data" rows for the PQA and FQA stages (mirroring the existing
# "Test" stage rows) on the "data" worksheet, and updates the active selections.

$wb = $excel.ActiveWorkbook

$dataSheet = $wb.Worksheets.Item("data")
$customerSheet = $wb.Worksheets.Item("CustomerData")

$hyperlinkAddress = "https://webdemo.baangt.org/"
$remoteReadValue = "RRD_[CustomerData,*,[Stage:[`$(Stage)],Country:[`$(Country)],IsActive:[X]]]"

# Each entry: row number, Stage value, then the C/E/F/G values that mirror the
# existing "Test" stage block (rows 2-4).
$newRows = @(
    @{ Row = 5;  Stage = "PQA"; C = "test area. Yeah!";  E = "this is an input";          F = "USA";    G = "Baangt is great!" },
    @{ Row = 6;  Stage = "PQA"; C = "Wonderful baangt!"; E = $null;                        F = "Canada"; G = "Baangt is fantastic!" },
    @{ Row = 7;  Stage = "PQA"; C = "This is perfect!";  E = "Another input, how cool!";   F = "USA";    G = "Baangt is the real deal!" },
    @{ Row = 8;  Stage = "FQA"; C = "test area. Yeah!";  E = "this is an input";          F = "USA";    G = "Baangt is great!" },
    @{ Row = 9;  Stage = "FQA"; C = "Wonderful baangt!"; E = $null;                        F = "Canada"; G = "Baangt is fantastic!" },
    @{ Row = 10; Stage = "FQA"; C = "This is perfect!";  E = "Another input, how cool!";   F = "USA";    G = "Baangt is the real deal!" }
)

# Copy the style (Consolas font) used by the existing F2 "country" cell so the
# new F cells reuse the same cell format instead of creating new ones.
$countryStyleSource = $dataSheet.Range("F2")
$countryStyleSource.Copy()

foreach ($entry in $newRows) {
    $r = $entry.Row

    $dataSheet.Range("A$r").Value = $entry.Stage
    $dataSheet.Range("B$r").Value = $remoteReadValue
    $dataSheet.Range("C$r").Value = $entry.C

    $dCell = $dataSheet.Range("D$r")
    $h = $dataSheet.Hyperlinks.Add($dCell, $hyperlinkAddress)
    $dCell.Value = $hyperlinkAddress
    $dCell.Style = "Hyperlink"

    if ($entry.E -ne $null) {
        $dataSheet.Range("E$r").Value = $entry.E
    }

    $fCell = $dataSheet.Range("F$r")
    $fCell.PasteSpecial(-4122)
    $fCell.Value = $entry.F

    $dataSheet.Range("G$r").Value = $entry.G
}

$excel.CutCopyMode = 0

# Update active selections to match the saved state.
[void]$dataSheet.Range("B13").Select()
[void]$customerSheet.Range("A5").Select()

[void]$dataSheet.Activate()
